$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- Row 51: Game-Changer ---
$ws.Range("A51").Value = "Game-Changer"
$ws.Range("B51").Value = "David McAdams"
$ws.Range("C51").NumberFormat = "m/d/yy"
$ws.Range("C51").Value = 43927
$ws.Range("D51").NumberFormat = "m/d/yy"
$ws.Range("D51").Value = 43930
$ws.Range("E51").Value = "Game theory;economics;strategy"
$ws.Range("F51").Value = "Hard Copy"
$ws.Range("G51").Value = "238 Pages"

# --- Row 52: Operation Paperclip ---
$ws.Range("A52").Value = "Operation Paperclip"
$ws.Range("B52").Value = "Annie Jacobsen"
$ws.Range("C52").NumberFormat = "m/d/yy"
$ws.Range("C52").Value = 43927
$ws.Range("D52").NumberFormat = "m/d/yy"
$ws.Range("D52").Value = 43930
$ws.Range("E52").Value = "nazis;history;science;holocaust;classified operations"
$ws.Range("F52").Value = "Audio"
$ws.Range("G52").Value = "19 Hours 26 Mins"

# --- Update view state to match new extent ---
$ws.Range("A53").Select()
